# Database updated with One degree plan
$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheet "Degree": remove rows 3-5 (ACS+DB, ACS+NF, ACS), keep only ACS+2
# ----------------------------------------------------------------------
$wsDegree = $wb.Worksheets.Item("Degree")
$wsDegree.Rows("3:5").Delete() | Out-Null
$wsDegree.Activate() | Out-Null
$wsDegree.Range("A3:D5").Select() | Out-Null

# ----------------------------------------------------------------------
# Sheet "DegreeRequirement": remove rows 15-49 (requirements for the
# degrees that no longer exist), keep only rows for DegreeId 1
# ----------------------------------------------------------------------
$wsDegreeReq = $wb.Worksheets.Item("DegreeRequirement")
$wsDegreeReq.Rows("15:49").Delete() | Out-Null
$wsDegreeReq.Activate() | Out-Null
$wsDegreeReq.Range("C18").Select() | Out-Null

# ----------------------------------------------------------------------
# Sheet "DegreeplanTermRequirement": no row data changes, just update
# the view (selection) like the original commit.
# ----------------------------------------------------------------------
$wsDPTR = $wb.Worksheets.Item("DegreeplanTermRequirement")
$wsDPTR.Activate() | Out-Null
$wsDPTR.Range("B17").Select() | Out-Null

# ----------------------------------------------------------------------
# Sheet "DegreePlan": update StudentId references from S530473 to
# S531367 for both existing degree plan rows.
# ----------------------------------------------------------------------
$wsDegreePlan = $wb.Worksheets.Item("DegreePlan")
$wsDegreePlan.Range("C2").Value = "S531367"
$wsDegreePlan.Range("C3").Value = "S531367"
$wsDegreePlan.Activate() | Out-Null
$wsDegreePlan.Range("C3").Select() | Out-Null

# ----------------------------------------------------------------------
# Sheet "Student": replace the student record with the new student.
# ----------------------------------------------------------------------
$wsStudent = $wb.Worksheets.Item("Student")
$wsStudent.Range("A2").Value = 531367
$wsStudent.Range("B2").Value = "Sai Sirisha"
$wsStudent.Range("C2").Value = "Devineni"
$wsStudent.Range("D2").Value = "s531367"
$wsStudent.Range("E2").Value = 562438
$wsStudent.Activate() | Out-Null
$wsStudent.Range("E2").Select() | Out-Null

# ----------------------------------------------------------------------
# Sheet "StudentTerm": shift term labels by one term (drop Fall 2017,
# add Summer 2019), update StudentID and add a new Fall 2019 row.
# ----------------------------------------------------------------------
$wsStudentTerm = $wb.Worksheets.Item("StudentTerm")
$wsStudentTerm.Range("B2").Value = 531367
$wsStudentTerm.Range("D2").Value = "Spring 2018"
$wsStudentTerm.Range("B3").Value = 531367
$wsStudentTerm.Range("D3").Value = "Summer 2018"
$wsStudentTerm.Range("B4").Value = 531367
$wsStudentTerm.Range("D4").Value = "Fall 2018"
$wsStudentTerm.Range("B5").Value = 531367
$wsStudentTerm.Range("D5").Value = "Spring 2019"
$wsStudentTerm.Range("B6").Value = 531367
$wsStudentTerm.Range("D6").Value = "Summer 2019"
$wsStudentTerm.Range("A7").Value = 6
$wsStudentTerm.Range("B7").Value = 531367
$wsStudentTerm.Range("C7").Value = 6
$wsStudentTerm.Range("D7").Value = "Fall 2019"
$wsStudentTerm.Columns("A").ColumnWidth = 16.28515625
$wsStudentTerm.Columns("B").ColumnWidth = 19
$wsStudentTerm.Columns("C").ColumnWidth = 16.85546875
$wsStudentTerm.Columns("D").ColumnWidth = 19.7109375
$wsStudentTerm.Activate() | Out-Null
$wsStudentTerm.Range("D1").Select() | Out-Null

# ----------------------------------------------------------------------
# Restore the originally active sheet (DegreeRequirement, tab index 3)
# ----------------------------------------------------------------------
$wsDegreeReq.Activate() | Out-Null
